$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values are plain decimal numbers need to be forced to
# Text format first, otherwise Excel auto-converts them to numbers and
# the exact textual representation (trailing zeros, etc.) would be lost.
$textForceCells = @("D4", "D5", "D6", "D7", "D9", "D10", "D11", "D12", "D14", "D15", "D19", "D20", "D21", "D23", "D24", "D25", "D26", "D27", "D28", "D30", "D31", "D32", "D33", "D34", "D35", "D37", "D39", "D40", "D41", "D44", "D45", "D46", "D47", "D48", "D49", "D51")
foreach ($ref in $textForceCells) {
    $ws.Range($ref).NumberFormat = "@"
}

$ws.Range("D2").Value = '51.570.61'
$ws.Range("E2").Value = '  +0.68%  '
$ws.Range("D3").Value = '2.982.62'
$ws.Range("E3").Value = '  +2.30%  '
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.07%  '
$ws.Range("D5").Value = '382.55'
$ws.Range("E5").Value = '  +2.31%  '
$ws.Range("D6").Value = '103.65'
$ws.Range("E6").Value = '  -0.43%  '
$ws.Range("D7").Value = '0.547'
$ws.Range("E7").Value = '  +0.60%  '
$ws.Range("E8").Value = '  +0.01%  '
$ws.Range("D9").Value = '0.594'
$ws.Range("E9").Value = '  +0.58%  '
$ws.Range("D10").Value = '37.10'
$ws.Range("E10").Value = '  -0.07%  '
$ws.Range("D11").Value = '0.139'
$ws.Range("E11").Value = '  +0.07%  '
$ws.Range("D12").Value = '0.0848'
$ws.Range("E12").Value = '  +1.22%  '
$ws.Range("D13").Value = '3.450.51'
$ws.Range("D14").Value = '18.28'
$ws.Range("E14").Value = '  -0.47%  '
$ws.Range("D15").Value = '7.54'
$ws.Range("E15").Value = '  +1.97%  '
$ws.Range("D16").Value = '2.970.55'
$ws.Range("E16").Value = '  +1.87%  '
$ws.Range("E17").Value = '  +7.57%  '
$ws.Range("D18").Value = '51.491.05'
$ws.Range("E18").Value = '  +0.55%  '
$ws.Range("D19").Value = '3.26'
$ws.Range("E19").Value = '  -1.72%  '
$ws.Range("D20").Value = '7.42'
$ws.Range("E20").Value = '  +1.95%  '
$ws.Range("D21").Value = '12.82'
$ws.Range("E21").Value = '  -1.31%  '
$ws.Range("D22").Value = '0.0₃0962'
$ws.Range("E22").Value = '  +1.71%  '
$ws.Range("D23").Value = '69.18'
$ws.Range("E23").Value = '  +1.03%  '
$ws.Range("D24").Value = '263.50'
$ws.Range("E24").Value = '  +1.14%  '
$ws.Range("D25").Value = '2.89'
$ws.Range("E25").Value = '  +7.18%  '
$ws.Range("D26").Value = '8.24'
$ws.Range("E26").Value = '  +14.87%  '
$ws.Range("D27").Value = '7.65'
$ws.Range("E27").Value = '  +15.55%  '
$ws.Range("D28").Value = '0.117'
$ws.Range("E28").Value = '  +15.07%  '
$ws.Range("E29").Value = '  -0.80%  '
$ws.Range("D30").Value = '4.14'
$ws.Range("E30").Value = '  +0.17%  '
$ws.Range("B31").Value = 'Dai'
$ws.Range("C31").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D31").Value = '1.00'
$ws.Range("E31").Value = '  +0.03%  '
$ws.Range("B32").Value = 'EthereumClassic'
$ws.Range("C32").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D32").Value = '26.02'
$ws.Range("E32").Value = '  +0.89%  '
$ws.Range("D33").Value = '9.86'
$ws.Range("E33").Value = '  -0.50%  '
$ws.Range("D34").Value = '34.58'
$ws.Range("E34").Value = '  +0.09%  '
$ws.Range("D35").Value = '51.03'
$ws.Range("E35").Value = '  +0.13%  '
$ws.Range("E36").Value = '  -2.26%  '
$ws.Range("D37").Value = '0.0451'
$ws.Range("E37").Value = '  +6.42%  '
$ws.Range("E38").Value = '  -0.21%  '
$ws.Range("D39").Value = '2.99'
$ws.Range("E39").Value = '  -0.31%  '
$ws.Range("D40").Value = '16.94'
$ws.Range("E40").Value = '  -0.85%  '
$ws.Range("D41").Value = '2.58'
$ws.Range("E41").Value = '  -0.17%  '
$ws.Range("E42").Value = '  +2.39%  '
$ws.Range("E43").Value = '  -0.80%  '
$ws.Range("D44").Value = '122.67'
$ws.Range("E44").Value = '  +2.69%  '
$ws.Range("D45").Value = '21.60'
$ws.Range("E45").Value = '  -1.39%  '
$ws.Range("B46").Value = 'WEMIXToken'
$ws.Range("C46").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D46").Value = '2.05'
$ws.Range("E46").Value = '  -2.11%  '
$ws.Range("B47").Value = 'TheGraph'
$ws.Range("C47").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D47").Value = '0.275'
$ws.Range("E47").Value = '  +9.56%  '
$ws.Range("B48").Value = 'ApeXProtocol'
$ws.Range("C48").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D48").Value = '2.37'
$ws.Range("E48").Value = '  +2.66%  '
$ws.Range("B49").Value = 'NEARProtocol'
$ws.Range("C49").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D49").Value = '3.31'
$ws.Range("E49").Value = '  +4.02%  '
$ws.Range("D50").Value = '2.034.08'
$ws.Range("E50").Value = '  +0.45%  '
$ws.Range("D51").Value = '0.0332'
$ws.Range("E51").Value = '  +4.63%  '
